# @williamjreid Minor checks and commenting updates
#
# Restructure the workbook's sheet tabs:
#   Exercise 2 | Exercise 4 | Sheet3   ->   Exercise 2 | Exercise 3 | Exercise 4 | Exercise 4(2)
#
# "Exercise 4" keeps its original identity/position data (just shifts one tab to the
# right). "Sheet3" is removed and replaced by a brand new "Exercise 4(2)" sheet (placed
# right after "Exercise 4") that gets a fresh small data set + AVERAGE() summaries.
# A new, empty "Exercise 3" placeholder sheet is inserted right before "Exercise 4".

$wb = $excel.ActiveWorkbook

# --- remove the old placeholder sheet -------------------------------------------------
$sheet3 = $wb.Worksheets.Item("Sheet3")
$sheet3.Delete() | Out-Null

# --- add "Exercise 4(2)" right after "Exercise 4" and fill it in ----------------------
$exercise4 = $wb.Worksheets.Item("Exercise 4")
$ex42 = $wb.Worksheets.Add($null, $exercise4)
$ex42.Name = "Exercise 4(2)"

$ex42.Columns.Item(1).ColumnWidth = 17

$ex42.Range("A1").Value = "a280_279"
$ex42.Range("A2").Value = 3473.92
$ex42.Range("A3").Value = 6164.97
$ex42.Range("A4").Value = 5616.14
$ex42.Range("A6").Formula = "=AVERAGE(A2:A4)"

$ex42.Range("A9").Value = "a280_1395"
$ex42.Range("A10").Value = -162238.44
$ex42.Range("A11").Value = -205335.97
$ex42.Range("A13").Formula = "=AVERAGE(A10:A11)"

$ex42.Range("A16").Value = "a280_2790"
$ex42.Range("A17").Value = -620626.9
$ex42.Range("A18").Value = -444818.84
$ex42.Range("A20").Formula = "=AVERAGE(A17:A18)"

$ex42.Range("A23").Value = "fnl4461_4460"
$ex42.Range("A24").Value = -6927376.77
$ex42.Range("A25").Value = -6926291.15
$ex42.Range("A27").Formula = "=AVERAGE(A24:A25)"

$ex42.Range("A29").Value = "fnl4461_22300"
$ex42.Range("A30").Value = -86880402.02
$ex42.Range("A31").Value = -86539272.71
$ex42.Range("A33").Formula = "=AVERAGE(A30:A31)"

$ex42.Activate() | Out-Null
$ex42.Range("C25").Select() | Out-Null

# --- add the blank "Exercise 3" placeholder right before "Exercise 4" -----------------
$exercise4b = $wb.Worksheets.Item("Exercise 4")
$ex3 = $wb.Worksheets.Add($exercise4b)
$ex3.Name = "Exercise 3"

$ex3.Activate() | Out-Null
$ex3.Range("D23").Select() | Out-Null

# --- tidy "Exercise 4"'s leftover selection from when it used to be the active tab ----
$exercise4c = $wb.Worksheets.Item("Exercise 4")
$exercise4c.Activate() | Out-Null
$exercise4c.Range("A1").Select() | Out-Null

# --- leave "Exercise 3" as the active tab, matching the saved selection state ---------
$ex3again = $wb.Worksheets.Item("Exercise 3")
$ex3again.Activate() | Out-Null
